$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1163
$ws1.Range("F5").Value = 66
$ws1.Range("F7").Value = 821
$ws1.Range("F8").Value = 431
$ws1.Range("F10").Value = 2070
$ws1.Range("F12").Value = 243
$ws1.Range("F14").Value = 965
$ws1.Range("F15").Value = 136
$ws1.Range("F16").Value = 2093
$ws1.Range("F17").Value = 563
$ws1.Range("F18").Value = 9866
$ws1.Range("F19").Value = 945
$ws1.Range("F20").Value = 539
$ws1.Range("F21").Value = 96
$ws1.Range("F22").Value = 116
$ws1.Range("F24").Value = 248

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3
$ws2.Range("F10").Value = 135
$ws2.Range("F11").Value = 8
$ws2.Range("F16").Value = 2
$ws2.Range("F19").Value = 1
$ws2.Range("F24").Value = 18

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5642
$ws3.Range("F3").Value = 450
$ws3.Range("F4").Value = 426

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5642
$ws4.Range("F4").Value = 450
$ws4.Range("F5").Value = 426
$ws4.Range("F7").Value = 1163
$ws4.Range("F9").Value = 3
$ws4.Range("F10").Value = 66
$ws4.Range("F12").Value = 821
$ws4.Range("F14").Value = 431
$ws4.Range("F16").Value = 2070
$ws4.Range("F18").Value = 243
$ws4.Range("F22").Value = 965
$ws4.Range("F24").Value = 136
$ws4.Range("F25").Value = 135
$ws4.Range("F26").Value = 8
$ws4.Range("F27").Value = 2093
$ws4.Range("F28").Value = 563
$ws4.Range("F31").Value = 945
$ws4.Range("F32").Value = 539
$ws4.Range("F33").Value = 96
$ws4.Range("F34").Value = 116
$ws4.Range("F38").Value = 2
$ws4.Range("F39").Value = 248
$ws4.Range("F42").Value = 1
$ws4.Range("F48").Value = 18
